$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("A2").Value = 'Peru'
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '1'
$ws.Range("C2").Value = 'Insurance (General)'
$ws.Range("D2").Value = 0.0632
$ws.Range("E2").Value = 0.05860000000000001
$ws.Range("G2").Value = 0.1339950372208437
$ws.Range("H2").Value = 0.1339950372208437
$ws.Range("I2").Value = 0.1153846153846154
$ws.Range("J2").Value = 0.08179868458678867
$ws.Range("K2").Value = 19
$ws.Range("L2").Value = 0.07857733664185276
$ws.Range("M2").Value = 23.5
$ws.Range("N2").Value = 0.2800953516090584
$ws.Range("O2").Value = 1.236842105263158
$ws.Range("P2").Value = 23.5
$ws.Range("Q2").Value = 0.2800953516090584
$ws.Range("R2").Value = 1.236842105263158
$ws.Range("U2").Value = 10.6
$ws.Range("V2").Value = 0.1263408820023838
$ws.Range("W2").Value = 0.1556101556101556
$ws.Range("X2").Value = 0.04817682442461023
$ws.Range("Y2").Value = 0.1074333311855454
$ws.Range("Z2").Value = 2.522823308466795
$ws.Range("AA2").Value = 0.2063636280774741
$ws.Range("AB2").Value = 0.04817682442461023
$ws.Range("AC2").Value = 0.1581868036528639
$ws.Range("AD2").Value = 0
$ws.Range("AF2").Value = 0
$ws.Range("AG2").Value = -10.6
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = -0.1446111869031378
$ws.Range("AK2").Value = -0.09532374100719423
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0
$ws.Range("AP2").Value = -0.3271604938271605

# --- Row 3 updates ---
$ws.Range("B3").Value = 'Mapfre Perú Compañía de Seguros y Reaseguros S.A. (BVL:MAPFREC1)'
$ws.Range("D3").Value = 0.0632
$ws.Range("E3").Value = 0.05860000000000001
$ws.Range("G3").Value = 0.1339950372208437
$ws.Range("H3").Value = 0.1339950372208437
$ws.Range("I3").Value = 0.1153846153846154
$ws.Range("J3").Value = 0.08179868458678867
$ws.Range("K3").Value = 19
$ws.Range("L3").Value = 0.07857733664185276
$ws.Range("M3").Value = 23.5
$ws.Range("N3").Value = 0.2800953516090584
$ws.Range("O3").Value = 1.236842105263158
$ws.Range("P3").Value = 23.5
$ws.Range("Q3").Value = 0.2800953516090584
$ws.Range("R3").Value = 1.236842105263158
$ws.Range("U3").Value = 10.6
$ws.Range("V3").Value = 0.1263408820023838
$ws.Range("W3").Value = 0.1556101556101556
$ws.Range("X3").Value = 0.04817682442461023
$ws.Range("Y3").Value = 0.1074333311855454
$ws.Range("Z3").Value = 2.522823308466795
$ws.Range("AA3").Value = 0.2063636280774741
$ws.Range("AB3").Value = 0.04817682442461023
$ws.Range("AC3").Value = 0.1581868036528639
$ws.Range("AD3").Value = 0
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = -10.6
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = -0.1446111869031378
$ws.Range("AK3").Value = -0.09532374100719423
$ws.Range("AL3").Value = 0
$ws.Range("AM3").Value = 0
$ws.Range("AN3").Value = 0
$ws.Range("AP3").Value = -0.3271604938271605

# --- Clear AO/AQ for rows 2 and 3 (columns removed) ---
$ws.Range("AO2").ClearContents()
$ws.Range("AQ2").ClearContents()
$ws.Range("AO3").ClearContents()
$ws.Range("AQ3").ClearContents()

# --- Delete row 4 entirely (Mapfre row merged/removed) ---
$ws.Rows.Item(4).Delete()
